$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compounds")

# Update the NCI Thesaurus "source_version" value (column E) from "25.07d" to "25.08d".
$ws.Range("E3").Value = "25.08d"
